$p = $ppt.ActivePresentation
$v = $p.HasNotesMaster
Write-Output ("type=" + $v.GetType().Name)
Write-Output ("val=[" + $v + "]")
if ($v) { Write-Output "truthy" } else { Write-Output "falsy" }
